# Aula 12 - Testando os links do menu - primeiros endpoints para acesso ao
# cadastro e lista de funcionario, cargo e departamento

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize the "sessao" column: every row belongs to session
# "2. Iniciando o desenvolvimento do projeto" (fixes a typo'd duplicate that
# existed in rows 4-6: "2. Inciando o Desenvolvimento do Projeto").
$ws.Range("C3").Value = "2. Iniciando o desenvolvimento do projeto"
$ws.Range("C4").Value = "2. Iniciando o desenvolvimento do projeto"
$ws.Range("C5").Value = "2. Iniciando o desenvolvimento do projeto"
$ws.Range("C6").Value = "2. Iniciando o desenvolvimento do projeto"

# Row 6 (aula 11) actually refers to the WebJars lesson.
$ws.Range("D6").Value = "11. Incluíndo o recurso de WebJars"
$ws.Range("E6").Value = "6:08 - foi ensinado uma forma de adicionar bibliotecas para pagina HTML (Jquery, bootstrap, icones) através de dependências no pom.xml. Usando esta forma, o `"src`" do documento html deve direcionar para o diretorio raiz onde encontra-se os webjars baixados pelo pom. É possivel ver o local destes diretórios direto na documentação das bibliotecas. Os webjars podem ser baixados atraves do site https://www.webjars.org/"

# New row 7: aula 12 - Testando os links do menu.
$ws.Range("B7").Value = 12
$ws.Range("C7").Value = "2. Iniciando o desenvolvimento do projeto"
$ws.Range("D7").Value = "Testando os links do menu"

$e7 = $ws.Range("E7")
$e7.Font.Underline = -4142
$e7.WrapText = $true
$e7.Value = "Nesta aula foi criado os primeiros endpoints para acesso dos menus de cadastro e lista de funcionario, cargos e departamentos"

# Match the row height used for the other multi-line, wrapped rows.
$ws.Rows.Item(7).RowHeight = 30

# Leave the selection where the author left it when saving.
$ws.Range("E11").Select()
